{"js": "// Update the \"Unidad\" title line from \"Unidad 1: Algoritmos\" to\n// \"Unidad 6: Creaci\u00f3n de aplicaciones de consola en java\".\nconst body = context.document.body;\n\n// The heading run \"Unidad 1:\" is bold; the following run \" Algoritmos\" is\n// regular. Find each literal run's text and replace it in place so the\n// bold/non-bold run split (and all other paragraph/run formatting) is kept.\nconst hits1 = body.search(\"Unidad 1:\", { matchCase: true });\nhits1.load(\"items\");\nconst hits2 = body.search(\" Algoritmos\", { matchCase: true });\nhits2.load(\"items\");\nawait context.sync();\n\nif (hits1.items.length === 0) {\n  throw new Error(\"Could not find 'Unidad 1:' text to update.\");\n}\nif (hits2.items.length === 0) {\n  throw new Error(\"Could not find ' Algoritmos' text to update.\");\n}\n\nhits1.items[0].insertText(\"Unidad 6:\", \"Replace\");\nhits2.items[0].insertText(\" Creaci\u00f3n de aplicaciones de consola en java\", \"Replace\");\n\nawait context.sync();\n", "ps1": "# Update the \"Unidad\" title line from \"Unidad 1: Algoritmos\" to\n# \"Unidad 6: Creaci\u00f3n de aplicaciones de consola en java\".\n#\n# The paragraph holds two runs:\n#   1. \"Unidad 1:\"   (bold)      -> \"Unidad 6:\"\n#   2. \" Algoritmos\" (regular)   -> \" Creaci\u00f3n de aplicaciones de consola en java\"\n# Find & Replace (scoped to exact text) edits each run's text in place, so\n# the existing bold / non-bold run split and all other formatting survive.\n\n$d = $word.ActiveDocument\n\n$find1 = $d.Content.Find\n$find1.Text = \"Unidad 1:\"\n$find1.Replacement.Text = \"Unidad 6:\"\n$find1.Execute(\"Unidad 1:\", $false, $false, $false, $false, $false, $true, 1, $false, \"Unidad 6:\", 2)\n\n$find2 = $d.Content.Find\n$find2.Text = \" Algoritmos\"\n$find2.Replacement.Text = \" Creaci\u00f3n de aplicaciones de consola en java\"\n$find2.Execute(\" Algoritmos\", $false, $false, $false, $false, $false, $true, 1, $false, \" Creaci\u00f3n de aplicaciones de consola en java\", 2)\n"}
